$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> FAPs
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Gnas"
$ws.Range("C2").Value = "Gcgr"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 136.287657
$ws.Range("H2").Value = 408.862971
$ws.Range("I2").Value = 0.2628768458810872
$ws.Range("J2").Value = 0.2628768458810872
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.8229729999999998
$ws.Range("N2").Value = 2.468919
$ws.Range("O2").Value = 0.887310079224403
$ws.Range("P2").Value = 0.8873100792244031
$ws.Range("Q2").Value = 112.161061944261
$ws.Range("R2").Value = 1009.449557498349
$ws.Range("S2").Value = 0.2332532749450087
$ws.Range("T2").Value = 0.2332532749450087

# Row 3: ECs -> sCs
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Gnas"
$ws.Range("C3").Value = "Gcgr"
$ws.Range("D3").Value = "sCs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 136.287657
$ws.Range("H3").Value = 408.862971
$ws.Range("I3").Value = 0.2628768458810872
$ws.Range("J3").Value = 0.2628768458810872
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.104519
$ws.Range("N3").Value = 0.313557
$ws.Range("O3").Value = 0.112689920775597
$ws.Range("P3").Value = 0.112689920775597
$ws.Range("Q3").Value = 14.244649621983
$ws.Range("R3").Value = 128.201846597847
$ws.Range("S3").Value = 0.02962357093607854
$ws.Range("T3").Value = 0.02962357093607854

# Row 4: FAPs -> FAPs
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Gnas"
$ws.Range("C4").Value = "Gcgr"
$ws.Range("D4").Value = "FAPs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 143.539174
$ws.Range("H4").Value = 430.617522
$ws.Range("I4").Value = 0.2768638492442244
$ws.Range("J4").Value = 0.2768638492442244
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.8229729999999998
$ws.Range("N4").Value = 2.468919
$ws.Range("O4").Value = 0.887310079224403
$ws.Range("P4").Value = 0.8873100792244031
$ws.Range("Q4").Value = 118.128864644302
$ws.Range("R4").Value = 1063.159781798718
$ws.Range("S4").Value = 0.2456640840072659
$ws.Range("T4").Value = 0.2456640840072659

# Row 5: FAPs -> sCs
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Gnas"
$ws.Range("C5").Value = "Gcgr"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 143.539174
$ws.Range("H5").Value = 430.617522
$ws.Range("I5").Value = 0.2768638492442244
$ws.Range("J5").Value = 0.2768638492442244
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.104519
$ws.Range("N5").Value = 0.313557
$ws.Range("O5").Value = 0.112689920775597
$ws.Range("P5").Value = 0.112689920775597
$ws.Range("Q5").Value = 15.002570927306
$ws.Range("R5").Value = 135.023138345754
$ws.Range("S5").Value = 0.03119976523695848
$ws.Range("T5").Value = 0.03119976523695847

# Row 6: sCs -> FAPs
$ws.Range("A6").Value = "sCs"
$ws.Range("B6").Value = "Gnas"
$ws.Range("C6").Value = "Gcgr"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 238.6199593333333
$ws.Range("H6").Value = 715.859878
$ws.Range("I6").Value = 0.4602593048746885
$ws.Range("J6").Value = 0.4602593048746884
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.8229729999999998
$ws.Range("N6").Value = 2.468919
$ws.Range("O6").Value = 0.887310079224403
$ws.Range("P6").Value = 0.8873100792244031
$ws.Range("Q6").Value = 196.3777837924313
$ws.Range("R6").Value = 1767.400054131882
$ws.Range("S6").Value = 0.4083927202721285
$ws.Range("T6").Value = 0.4083927202721285

# Row 7: sCs -> sCs
$ws.Range("A7").Value = "sCs"
$ws.Range("B7").Value = "Gnas"
$ws.Range("C7").Value = "Gcgr"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 238.6199593333333
$ws.Range("H7").Value = 715.859878
$ws.Range("I7").Value = 0.4602593048746885
$ws.Range("J7").Value = 0.4602593048746884
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.3333333333333333
$ws.Range("M7").Value = 0.104519
$ws.Range("N7").Value = 0.313557
$ws.Range("O7").Value = 0.112689920775597
$ws.Range("P7").Value = 0.112689920775597
$ws.Range("Q7").Value = 24.94031952956066
$ws.Range("R7").Value = 224.462875766046
$ws.Range("S7").Value = 0.05186658460255998
$ws.Range("T7").Value = 0.05186658460255998

